# LOM3246.xlsx - reorder the "Requisitos" list so that the
# "LOM3229 - Metodos Experimentais da Fisica II (Indicacao de Conjunto)"
# entry moves to the top, with the other two requisitos shifting down.
#
# Before:
#   B23/C23 = LOB1021 -  Física IV  (Requisito)
#   B24/C24 = LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)
#   B25/C25 = LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)
#
# After:
#   B23/C23 = LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)
#   B24/C24 = LOB1021 -  Física IV  (Requisito)
#   B25/C25 = LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$line1 = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"
$line2 = "LOB1021 -  Física IV  (Requisito)`n"
$line3 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $line1
$ws.Range("C23").Value = $line1

$ws.Range("B24").Value = $line2
$ws.Range("C24").Value = $line2

$ws.Range("B25").Value = $line3
$ws.Range("C25").Value = $line3
